$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.540099999999993
$ws.Range("A8").Value = -21.06740000000001
$ws.Range("A10").Value = -20.57949999999997
$ws.Range("A12").Value = -22.38040000000003
$ws.Range("C13").Value = -12.79869999999998
$ws.Range("A18").Value = -22.20110000000002
$ws.Range("D20").Value = -8.413899999999995
$ws.Range("A25").Value = -22.29000000000003
